$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds a date serial number that was bumped
# by one day (45203 -> 45204) for every data row (2 through 427).
$ws.Range("C2:C427").Value = 45204
